$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.206015333333333
$ws.Range("H2").Value = 18.618046
$ws.Range("I2").Value = 0.0150172404156507
$ws.Range("J2").Value = 0.0150172404156507
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2615913333333333
$ws.Range("N2").Value = 0.784774
$ws.Range("O2").Value = 0.08239613548481725
$ws.Range("P2").Value = 0.08239613548481727
$ws.Range("Q2").Value = 1.623439825733778
$ws.Range("R2").Value = 14.610958431604
$ws.Range("S2").Value = 0.001237362575896028
$ws.Range("T2").Value = 0.001237362575896028
$ws.Range("G3").Value = 6.206015333333333
$ws.Range("H3").Value = 18.618046
$ws.Range("I3").Value = 0.0150172404156507
$ws.Range("J3").Value = 0.0150172404156507
$ws.Range("N3").Value = 5.233242000000001
$ws.Range("O3").Value = 0.5494561706387266
$ws.Range("P3").Value = 0.5494561706387268
$ws.Range("Q3").Value = 10.82586003168133
$ws.Range("R3").Value = 97.43274028513201
$ws.Range("S3").Value = 0.008251315412344553
$ws.Range("T3").Value = 0.008251315412344553
$ws.Range("G4").Value = 6.206015333333333
$ws.Range("H4").Value = 18.618046
$ws.Range("I4").Value = 0.0150172404156507
$ws.Range("J4").Value = 0.0150172404156507
$ws.Range("M4").Value = 1.168795666666667
$ws.Range("N4").Value = 3.506387
$ws.Range("O4").Value = 0.3681476938764561
$ws.Range("P4").Value = 0.3681476938764561
$ws.Range("Q4").Value = 7.25356382886689
$ws.Range("R4").Value = 65.282074459802
$ws.Range("S4").Value = 0.005528562427410118
$ws.Range("T4").Value = 0.005528562427410118
$ws.Range("I5").Value = 0.9317452840597572
$ws.Range("J5").Value = 0.9317452840597571
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2615913333333333
$ws.Range("N5").Value = 0.784774
$ws.Range("O5").Value = 0.08239613548481725
$ws.Range("P5").Value = 0.08239613548481727
$ws.Range("Q5").Value = 100.7263891177905
$ws.Range("R5").Value = 906.5375020601141
$ws.Range("S5").Value = 0.0767722106627273
$ws.Range("T5").Value = 0.0767722106627273
$ws.Range("I6").Value = 0.9317452840597572
$ws.Range("J6").Value = 0.9317452840597571
$ws.Range("N6").Value = 5.233242000000001
$ws.Range("O6").Value = 0.5494561706387266
$ws.Range("P6").Value = 0.5494561706387268
$ws.Range("Q6").Value = 671.6909199840514
$ws.Range("R6").Value = 6045.218279856464
$ws.Range("S6").Value = 0.5119531957901667
$ws.Range("T6").Value = 0.5119531957901668
$ws.Range("I7").Value = 0.9317452840597572
$ws.Range("J7").Value = 0.9317452840597571
$ws.Range("M7").Value = 1.168795666666667
$ws.Range("N7").Value = 3.506387
$ws.Range("O7").Value = 0.3681476938764561
$ws.Range("P7").Value = 0.3681476938764561
$ws.Range("Q7").Value = 450.0476587648953
$ws.Range("R7").Value = 4050.428928884058
$ws.Range("S7").Value = 0.3430198776068631
$ws.Range("T7").Value = 0.3430198776068631
$ws.Range("G8").Value = 22.00088566666667
$ws.Range("H8").Value = 66.002657
$ws.Range("I8").Value = 0.05323747552459213
$ws.Range("J8").Value = 0.05323747552459213
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2615913333333333
$ws.Range("N8").Value = 0.784774
$ws.Range("O8").Value = 0.08239613548481725
$ws.Range("P8").Value = 0.08239613548481727
$ws.Range("Q8").Value = 5.755241016057555
$ws.Range("R8").Value = 51.79716914451799
$ws.Range("S8").Value = 0.004386562246193936
$ws.Range("T8").Value = 0.004386562246193936
$ws.Range("G9").Value = 22.00088566666667
$ws.Range("H9").Value = 66.002657
$ws.Range("I9").Value = 0.05323747552459213
$ws.Range("J9").Value = 0.05323747552459213
$ws.Range("N9").Value = 5.233242000000001
$ws.Range("O9").Value = 0.5494561706387266
$ws.Range("P9").Value = 0.5494561706387268
$ws.Range("Q9").Value = 38.37865296933267
$ws.Range("R9").Value = 345.407876723994
$ws.Range("S9").Value = 0.02925165943621533
$ws.Range("T9").Value = 0.02925165943621533
$ws.Range("G10").Value = 22.00088566666667
$ws.Range("H10").Value = 66.002657
$ws.Range("I10").Value = 0.05323747552459213
$ws.Range("J10").Value = 0.05323747552459213
$ws.Range("M10").Value = 1.168795666666667
$ws.Range("N10").Value = 3.506387
$ws.Range("O10").Value = 0.3681476938764561
$ws.Range("P10").Value = 0.3681476938764561
$ws.Range("Q10").Value = 25.71453983002878
$ws.Range("R10").Value = 231.430858470259
$ws.Range("S10").Value = 0.01959925384218287
$ws.Range("T10").Value = 0.01959925384218287
